$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row containing "syntok" in column A and delete the entire row.
$target = $ws.Range("A1:A38").Find("syntok")
if ($target -ne $null) {
    $target.EntireRow.Delete()
}

# Update selection / view to match the final state.
$ws.Application.GoTo($ws.Range("D34"), $true)
$ws.Range("D34").Select()
